# "Update countries & provincias Spain"
# Refreshes the COVID country statistics table on sheet "Pais":
#  - bumps the "Datos actualizados ..." timestamp in A1
#  - updates the numeric columns (B..H) for the countries whose figures changed
#  - because the table is kept sorted descending by column B (Casos totales),
#    two pairs of countries (Albania/Paraguay and Namibia/Sudan del Sur) swap
#    row order; we simply write the new country name + stats into each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Agosto de 2020 a las 16:41"

# Row 4
$ws.Range("B4").Value = 4865523
$ws.Range("C4").Value = 3349
$ws.Range("D4").Value = 2448515
$ws.Range("E4").Value = 2257992
$ws.Range("G4").Value = 88
$ws.Range("H4").Value = 159016

# Row 46
$ws.Range("D46").Value = 47454
$ws.Range("E46").Value = 5865

# Row 59
$ws.Range("B59").Value = 32910
$ws.Range("C59").Value = 226
$ws.Range("D59").Value = 28348
$ws.Range("E59").Value = 4089
$ws.Range("G59").Value = 5
$ws.Range("H59").Value = 473

# Row 66
$ws.Range("B66").Value = 23202
$ws.Range("C66").Value = 605
$ws.Range("D66").Value = 9327
$ws.Range("E66").Value = 13487
$ws.Range("G66").Value = 6
$ws.Range("H66").Value = 388

# Row 86
$ws.Range("B86").Value = 9348
$ws.Range("C86").Value = 14
$ws.Range("E86").Value = 340

# Row 91
$ws.Range("B91").Value = 7583
$ws.Range("C91").Value = 45
$ws.Range("D91").Value = 6356
$ws.Range("E91").Value = 1166

# Row 96
$ws.Range("B96").Value = 6793
$ws.Range("C96").Value = 213
$ws.Range("D96").Value = 5109
$ws.Range("E96").Value = 1511
$ws.Range("G96").Value = 2
$ws.Range("H96").Value = 173

# Row 98
$ws.Range("A98").Value = "Albania"
$ws.Range("B98").Value = 5750
$ws.Range("C98").Value = 130
$ws.Range("D98").Value = 3031
$ws.Range("E98").Value = 2543
$ws.Range("G98").Value = 4
$ws.Range("H98").Value = 176

# Row 99
$ws.Range("A99").Value = "Paraguay"
$ws.Range("B99").Value = 5724
$ws.Range("D99").Value = 4249
$ws.Range("E99").Value = 1420
$ws.Range("H99").Value = 55

# Row 101
$ws.Range("B101").Value = 5248
$ws.Range("C101").Value = 8
$ws.Range("D101").Value = 5044
$ws.Range("E101").Value = 145

# Row 123
$ws.Range("A123").Value = "Namibia"
$ws.Range("B123").Value = 2470
$ws.Range("C123").Value = 64
$ws.Range("D123").Value = 211
$ws.Range("E123").Value = 2247
$ws.Range("H123").Value = 12

# Row 124
$ws.Range("A124").Value = "Sudan del Sur"
$ws.Range("B124").Value = 2429
$ws.Range("D124").Value = 1175
$ws.Range("E124").Value = 1208
$ws.Range("H124").Value = 46

# Row 150
$ws.Range("B150").Value = 939
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 825
$ws.Range("E150").Value = 62
